$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Sheet: Recommandations ---
$ws1.Range("C2").Value = 4
$ws1.Range("D2").Value = 957.45
$ws1.Range("A3").Value = 'BRVM - CONSOMMATION DISCRETIONNAIRE'
$ws1.Range("D3").Value = 729.64
$ws1.Range("E3").Value = 188.51
$ws1.Range("A4").Value = 'BRVM - CONSOMMATION DE BASE     (**)'
$ws1.Range("C4").Value = 3
$ws1.Range("D4").Value = 719
$ws1.Range("E4").Value = 243.54
$ws1.Range("C5").Value = 4
$ws1.Range("D5").Value = 680.13
$ws1.Range("C6").Value = 4
$ws1.Range("D6").Value = 625.6
$ws1.Range("C7").Value = 4
$ws1.Range("D7").Value = 605.66
$ws1.Range("C8").Value = 4
$ws1.Range("D8").Value = 483.43
$ws1.Range("C9").Value = 4
$ws1.Range("D9").Value = 470.09
$ws1.Range("C10").Value = 3
$ws1.Range("D10").Value = 427.4
$ws1.Range("C11").Value = 4
$ws1.Range("D11").Value = 393.48
$ws1.Range("B14").Value = 4
$ws1.Range("D14").Value = 26.58
$ws1.Range("A18").Value = 'SICOR CI (SICC)'
$ws1.Range("B18").Value = 2
$ws1.Range("D18").Value = 13.85
$ws1.Range("E18").Value = 7.47
$ws1.Range("A19").Value = 'ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)'
$ws1.Range("C19").Value = 0
$ws1.Range("D19").Value = 9.880000000000001
$ws1.Range("E19").Value = 7.44
$ws1.Range("G19").Value = '➖ Neutre'
$ws1.Range("A20").Value = 'SETAO CI (STAC)'
$ws1.Range("B20").Value = 1
$ws1.Range("C20").Value = 0
$ws1.Range("D20").Value = 7.25
$ws1.Range("E20").Value = 7.25
$ws1.Range("G20").Value = '➖ Neutre'
$ws1.Range("A21").Value = 'SAFCA CI (SAFC)'
$ws1.Range("C21").Value = 0
$ws1.Range("D21").Value = 2.71
$ws1.Range("E21").Value = 2.71
$ws1.Range("G21").Value = '➖ Neutre'
$ws1.Range("A22").Value = 'FILTISAC CI (FTSC)'
$ws1.Range("C22").Value = 0
$ws1.Range("D22").Value = 2.62
$ws1.Range("E22").Value = 2.62
$ws1.Range("G22").Value = '➖ Neutre'
$ws1.Range("A23").Value = 'LOTERIE NATIONALE DU BENIN (LNBB)'
$ws1.Range("B23").Value = 0
$ws1.Range("D23").Value = -1.02
$ws1.Range("E23").Value = -1.02
$ws1.Range("G23").Value = '➖ Neutre'
$ws1.Range("A24").Value = 'BANK OF AFRICA NG (BOAN)'
$ws1.Range("B24").Value = 0
$ws1.Range("C24").Value = 1
$ws1.Range("D24").Value = -1.7
$ws1.Range("E24").Value = -1.7
$ws1.Range("G24").Value = '➖ Neutre'
$ws1.Range("A25").Value = 'CIE CI (CIEC)'
$ws1.Range("D25").Value = -1.86
$ws1.Range("E25").Value = -1.86
$ws1.Range("A26").Value = 'CORIS BANK INTERNATIONAL (CBIBF)'
$ws1.Range("D26").Value = -2.04
$ws1.Range("E26").Value = -2.04
$ws1.Range("A27").Value = 'ONATEL BF (ONTBF)'
$ws1.Range("D27").Value = -2.21
$ws1.Range("E27").Value = -2.21
$ws1.Range("A28").Value = 'SAPH CI (SPHC)'
$ws1.Range("D28").Value = -2.4
$ws1.Range("E28").Value = -2.4
$ws1.Range("A29").Value = 'BANK OF AFRICA ML (BOAM)'
$ws1.Range("D29").Value = -2.41
$ws1.Range("E29").Value = -2.41
$ws1.Range("A30").Value = 'SITAB CI (STBC)'
$ws1.Range("B30").Value = 1
$ws1.Range("D30").Value = -2.53
$ws1.Range("E30").Value = 4.97
$ws1.Range("G30").Value = '👀 À surveiller'
$ws1.Range("A31").Value = 'BERNABE CI (BNBC)'
$ws1.Range("D31").Value = -3.25
$ws1.Range("E31").Value = -3.25
$ws1.Range("A32").Value = 'AFRICA GLOBAL LOGISTICS CI (SDSC)'
$ws1.Range("B32").Value = 0
$ws1.Range("D32").Value = -3.34
$ws1.Range("E32").Value = -3.34
$ws1.Range("G32").Value = '➖ Neutre'
$ws1.Range("A33").Value = 'BANK OF AFRICA BF (BOABF)'
$ws1.Range("C33").Value = 2
$ws1.Range("D33").Value = -3.99
$ws1.Range("E33").Value = -1.9
$ws1.Range("A34").Value = 'ECOBANK TRANS. INCORP. TG (ETIT)'
$ws1.Range("B34").Value = 0
$ws1.Range("C34").Value = 1
$ws1.Range("D34").Value = -4.17
$ws1.Range("E34").Value = -4.17
$ws1.Range("G34").Value = '➖ Neutre'
$ws1.Range("A35").Value = 'SODE CI (SDCC)'
$ws1.Range("B35").Value = 1
$ws1.Range("C35").Value = 2
$ws1.Range("D35").Value = -4.3
$ws1.Range("E35").Value = -6.89
$ws1.Range("G35").Value = '👀 À surveiller'
$ws1.Range("A36").Value = 'SOGB CI (SOGC)'
$ws1.Range("C36").Value = 1
$ws1.Range("D36").Value = -4.34
$ws1.Range("E36").Value = -4.34
$ws1.Range("A37").Value = 'ORAGROUP TOGO (ORGT)'
$ws1.Range("B37").Value = 0
$ws1.Range("D37").Value = -4.61
$ws1.Range("E37").Value = -2.61
$ws1.Range("G37").Value = '➖ Neutre'
$ws1.Range("A38").Value = 'NEI-CEDA CI (NEIC)'
$ws1.Range("B38").Value = 0
$ws1.Range("D38").Value = -10.59
$ws1.Range("E38").Value = -7.49
$ws1.Range("G38").Value = '➖ Neutre'

# --- Sheet: Top_YTD ---
$ws2.Range("B2").Value = 13161.65
$ws2.Range("B3").Value = 6258.85
$ws2.Range("B4").Value = 5209.11
$ws2.Range("A5").Value = 'BRVM - SERVICES FINANCIERS'
$ws2.Range("B5").Value = 4221.77
$ws2.Range("A6").Value = 'BRVM-PRESTIGE'
$ws2.Range("B6").Value = 3895.35
$ws2.Range("A7").Value = 'BRVM - CONSOMMATION DE BASE     (**)'
$ws2.Range("B7").Value = 3818.47
$ws2.Range("B8").Value = 2279.01
$ws2.Range("B9").Value = 2138.67
$ws2.Range("A10").Value = 'BRVM - TELECOMMUNICATIONS'
$ws2.Range("B10").Value = 1448.46
$ws2.Range("A11").Value = 'BRVM – COMPOSITE TOTAL RETURN     (**)'
$ws2.Range("B11").Value = 1325.44

